$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.516.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.04%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.689.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3893"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4020"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.48%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.488"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9976"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08768"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +10.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.520"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.97%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.053"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001347"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.672.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "97.98"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07236"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.280"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.490.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.033"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.73%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.344"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.48%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.66"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.648"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.359"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.61%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "138.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.856.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.08759"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.358"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.047"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.03032"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +10.44%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.995"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.61%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2770"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.84"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09141"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8018"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.31%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.477"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.49%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.632"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7246"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.96%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.283"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.401"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.47%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.9999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.28%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.08071"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.02%  "

